$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1963.7142
$ws.Cells.Item(2, 9).Value = 439.75
$ws.Cells.Item(2, 10).Value = 2901.5386
$ws.Cells.Item(2, 11).Value = 439.75
$ws.Cells.Item(2, 12).Value = 2901.5386
$ws.Cells.Item(2, 13).Value = -326.75
$ws.Cells.Item(2, 14).Value = -3127.5386

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 464.375
$ws.Cells.Item(9, 9).Value = 491.5
$ws.Cells.Item(9, 10).Value = 274.5
$ws.Cells.Item(9, 11).Value = 491.5
$ws.Cells.Item(9, 12).Value = 274.5
$ws.Cells.Item(9, 13).Value = -322.5
$ws.Cells.Item(9, 14).Value = -612.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1895.6666
$ws.Cells.Item(28, 9).Value = 1895.6666
$ws.Cells.Item(28, 11).Value = 1895.6666
$ws.Cells.Item(28, 13).Value = -1410.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 555.5714
$ws.Cells.Item(58, 9).Value = 77.8
$ws.Cells.Item(58, 10).Value = 1750
$ws.Cells.Item(58, 11).Value = 233.4
$ws.Cells.Item(58, 12).Value = 5250
$ws.Cells.Item(58, 13).Value = -83.39999999999998
$ws.Cells.Item(58, 14).Value = -5550

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4579.6
$ws.Cells.Item(64, 10).Value = 5499.6665
$ws.Cells.Item(64, 12).Value = 5499.6665
$ws.Cells.Item(64, 14).Value = -5995.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 4579.6
$ws.Cells.Item(67, 10).Value = 5499.6665
$ws.Cells.Item(67, 12).Value = 5499.6665
$ws.Cells.Item(67, 14).Value = -7215.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 148389.7
$ws.Cells.Item(74, 9).Value = 209783.58
$ws.Cells.Item(74, 10).Value = 5137.3335
$ws.Cells.Item(74, 11).Value = 209783.58
$ws.Cells.Item(74, 12).Value = 5137.3335
$ws.Cells.Item(74, 13).Value = -208847.58
$ws.Cells.Item(74, 14).Value = -7009.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 148389.7
$ws.Cells.Item(77, 9).Value = 209783.58
$ws.Cells.Item(77, 10).Value = 5137.3335
$ws.Cells.Item(77, 11).Value = 1048917.9
$ws.Cells.Item(77, 12).Value = 25686.6675
$ws.Cells.Item(77, 13).Value = -1044237.9
$ws.Cells.Item(77, 14).Value = -35046.6675

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 8980.4
$ws.Cells.Item(80, 9).Value = 6668
$ws.Cells.Item(80, 10).Value = 12449
$ws.Cells.Item(80, 11).Value = 20004
$ws.Cells.Item(80, 12).Value = 37347
$ws.Cells.Item(80, 13).Value = -19006
$ws.Cells.Item(80, 14).Value = -39343

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 8980.4
$ws.Cells.Item(83, 9).Value = 6668
$ws.Cells.Item(83, 10).Value = 12449
$ws.Cells.Item(83, 11).Value = 60012
$ws.Cells.Item(83, 12).Value = 112041
$ws.Cells.Item(83, 13).Value = -55020
$ws.Cells.Item(83, 14).Value = -122025

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3162.625
$ws.Cells.Item(98, 9).Value = 1175.5
$ws.Cells.Item(98, 10).Value = 3825
$ws.Cells.Item(98, 11).Value = 1175.5
$ws.Cells.Item(98, 12).Value = 3825
$ws.Cells.Item(98, 13).Value = 322.5
$ws.Cells.Item(98, 14).Value = -6821

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3696
$ws.Cells.Item(113, 9).Value = 3495
$ws.Cells.Item(113, 10).Value = 4500
$ws.Cells.Item(113, 11).Value = 3495
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = -241
$ws.Cells.Item(113, 14).Value = -11008

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 3162.625
$ws.Cells.Item(122, 9).Value = 1175.5
$ws.Cells.Item(122, 10).Value = 3825
$ws.Cells.Item(122, 11).Value = 3526.5
$ws.Cells.Item(122, 12).Value = 11475
$ws.Cells.Item(122, 13).Value = -1076.5
$ws.Cells.Item(122, 14).Value = -16375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 199
$ws.Cells.Item(12, 9).Value = 199
$ws.Cells.Item(12, 11).Value = 199
$ws.Cells.Item(12, 13).Value = -26

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 15798
$ws.Cells.Item(17, 9).Value = 12996.667
$ws.Cells.Item(17, 10).Value = 20000
$ws.Cells.Item(17, 11).Value = 12996.667
$ws.Cells.Item(17, 12).Value = 20000
$ws.Cells.Item(17, 13).Value = -12823.667
$ws.Cells.Item(17, 14).Value = -20346

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(82, 8).Value = 44999
$ws.Cells.Item(82, 10).Value = 44999
$ws.Cells.Item(82, 12).Value = 44999
$ws.Cells.Item(82, 14).Value = -45721

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(85, 8).Value = 44999
$ws.Cells.Item(85, 10).Value = 44999
$ws.Cells.Item(85, 12).Value = 44999
$ws.Cells.Item(85, 14).Value = -47495

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(135, 8).Value = 197498.67
$ws.Cells.Item(135, 10).Value = 197498.67
$ws.Cells.Item(135, 12).Value = 197498.67
$ws.Cells.Item(135, 14).Value = -207638.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 698.8333
$ws.Cells.Item(5, 9).Value = 658.8
$ws.Cells.Item(5, 10).Value = 899
$ws.Cells.Item(5, 11).Value = 658.8
$ws.Cells.Item(5, 12).Value = 899
$ws.Cells.Item(5, 13).Value = -545.8
$ws.Cells.Item(5, 14).Value = -1125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 601.55554
$ws.Cells.Item(22, 9).Value = 471.75
$ws.Cells.Item(22, 11).Value = 471.75
$ws.Cells.Item(22, 13).Value = -298.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(74, 8).Value = 48000
$ws.Cells.Item(74, 10).Value = 48000
$ws.Cells.Item(74, 12).Value = 48000
$ws.Cells.Item(74, 14).Value = -49872

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(77, 8).Value = 48000
$ws.Cells.Item(77, 10).Value = 48000
$ws.Cells.Item(77, 12).Value = 144000
$ws.Cells.Item(77, 14).Value = -153360

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 42954.8
$ws.Cells.Item(81, 10).Value = 42954.8
$ws.Cells.Item(81, 12).Value = 42954.8
$ws.Cells.Item(81, 14).Value = -45076.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(84, 8).Value = 42954.8
$ws.Cells.Item(84, 10).Value = 42954.8
$ws.Cells.Item(84, 12).Value = 128864.4
$ws.Cells.Item(84, 14).Value = -139472.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1702.3
$ws.Cells.Item(94, 9).Value = 1671.4445
$ws.Cells.Item(94, 11).Value = 1671.4445
$ws.Cells.Item(94, 13).Value = -1220.4445

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(122, 8).Value = 70258.336
$ws.Cells.Item(122, 10).Value = 70258.336
$ws.Cells.Item(122, 12).Value = 70258.336
$ws.Cells.Item(122, 14).Value = -80058.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 184404.25
$ws.Cells.Item(2, 9).Value = 220127.8
$ws.Cells.Item(2, 11).Value = 1320766.8
$ws.Cells.Item(2, 13).Value = -1320653.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 14.833333
$ws.Cells.Item(12, 9).Value = 4
$ws.Cells.Item(12, 10).Value = 17
$ws.Cells.Item(12, 11).Value = 12
$ws.Cells.Item(12, 12).Value = 51
$ws.Cells.Item(12, 13).Value = 161
$ws.Cells.Item(12, 14).Value = -397

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 1038.6923
$ws.Cells.Item(14, 9).Value = 1038.6923
$ws.Cells.Item(14, 11).Value = 3116.0769
$ws.Cells.Item(14, 13).Value = -2943.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 102.5
$ws.Cells.Item(33, 9).Value = 70
$ws.Cells.Item(33, 10).Value = 135
$ws.Cells.Item(33, 11).Value = 420
$ws.Cells.Item(33, 12).Value = 810
$ws.Cells.Item(33, 13).Value = -137
$ws.Cells.Item(33, 14).Value = -1376

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 11727
$ws.Cells.Item(117, 10).Value = 25696.75
$ws.Cells.Item(117, 12).Value = 77090.25
$ws.Cells.Item(117, 14).Value = -83974.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2870.5833
$ws.Cells.Item(131, 10).Value = 2870.5833
$ws.Cells.Item(131, 12).Value = 8611.749899999999
$ws.Cells.Item(131, 14).Value = -18691.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 9049.556
$ws.Cells.Item(134, 9).Value = 1655
$ws.Cells.Item(134, 10).Value = 18292.75
$ws.Cells.Item(134, 11).Value = 4965
$ws.Cells.Item(134, 12).Value = 54878.25
$ws.Cells.Item(134, 13).Value = 105
$ws.Cells.Item(134, 14).Value = -65018.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 1274
$ws.Cells.Item(139, 9).Value = 1274
$ws.Cells.Item(139, 11).Value = 3822
$ws.Cells.Item(139, 13).Value = 1318

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 5396.55
$ws.Cells.Item(140, 9).Value = 1281
$ws.Cells.Item(140, 11).Value = 3843
$ws.Cells.Item(140, 13).Value = 1337

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 61.636364
$ws.Cells.Item(2, 9).Value = 85.85714
$ws.Cells.Item(2, 11).Value = 85.85714
$ws.Cells.Item(2, 13).Value = 27.14286

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 32015
$ws.Cells.Item(47, 10).Value = 32015
$ws.Cells.Item(47, 12).Value = 32015
$ws.Cells.Item(47, 14).Value = -33151

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 32601.143
$ws.Cells.Item(62, 9).Value = 34552
$ws.Cells.Item(62, 11).Value = 34552
$ws.Cells.Item(62, 13).Value = -33866

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(65, 8).Value = 32601.143
$ws.Cells.Item(65, 9).Value = 34552
$ws.Cells.Item(65, 11).Value = 103656
$ws.Cells.Item(65, 13).Value = -100224

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1233.8
$ws.Cells.Item(102, 9).Value = 1282.5555
$ws.Cells.Item(102, 11).Value = 1282.5555
$ws.Cells.Item(102, 13).Value = 339.4445000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1774.5
$ws.Cells.Item(22, 9).Value = 1700
$ws.Cells.Item(22, 10).Value = 1799.3334
$ws.Cells.Item(22, 11).Value = 1700
$ws.Cells.Item(22, 12).Value = 1799.3334
$ws.Cells.Item(22, 13).Value = -1405
$ws.Cells.Item(22, 14).Value = -2389.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 1774.5
$ws.Cells.Item(27, 9).Value = 1700
$ws.Cells.Item(27, 10).Value = 1799.3334
$ws.Cells.Item(27, 11).Value = 1700
$ws.Cells.Item(27, 12).Value = 1799.3334
$ws.Cells.Item(27, 13).Value = -1593
$ws.Cells.Item(27, 14).Value = -2013.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4201.375
$ws.Cells.Item(40, 9).Value = 2411
$ws.Cells.Item(40, 10).Value = 4798.1665
$ws.Cells.Item(40, 11).Value = 2411
$ws.Cells.Item(40, 12).Value = 4798.1665
$ws.Cells.Item(40, 13).Value = -2275
$ws.Cells.Item(40, 14).Value = -5070.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 52157.2
$ws.Cells.Item(46, 9).Value = 73510.28999999999
$ws.Cells.Item(46, 10).Value = 2333.3333
$ws.Cells.Item(46, 11).Value = 73510.28999999999
$ws.Cells.Item(46, 12).Value = 2333.3333
$ws.Cells.Item(46, 13).Value = -73322.28999999999
$ws.Cells.Item(46, 14).Value = -2709.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(63, 8).Value = 90077
$ws.Cells.Item(63, 9).Value = 90077
$ws.Cells.Item(63, 11).Value = 90077
$ws.Cells.Item(63, 13).Value = -89328

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(66, 8).Value = 90077
$ws.Cells.Item(66, 9).Value = 90077
$ws.Cells.Item(66, 11).Value = 270231
$ws.Cells.Item(66, 13).Value = -266487

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(110, 8).Value = 39999
$ws.Cells.Item(110, 10).Value = 39999
$ws.Cells.Item(110, 12).Value = 39999
$ws.Cells.Item(110, 14).Value = -48179

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 60000
$ws.Cells.Item(109, 10).Value = 60000
$ws.Cells.Item(109, 12).Value = 60000
$ws.Cells.Item(109, 14).Value = -62774

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 5088.8237
$ws.Cells.Item(126, 10).Value = 6186.5
$ws.Cells.Item(126, 12).Value = 18559.5
$ws.Cells.Item(126, 14).Value = -23499.5
